# Update tracking numbers in "CheetahProcessing" sheet (15th June 2022 changes)
# Column C = ShipmentTrackNum, Column D = PackageTrackNum
# For most rows only column C changes; for a subset of rows column D mirrors
# the same new value as column C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new ShipmentTrackNum (column C) value
$newTrackNums = [ordered]@{
    2  = "320018813081"
    3  = "320018813092"
    4  = "320018813129"
    5  = "320018813140"
    6  = "320018813184"
    7  = "320018813200"
    8  = "320018813232"
    9  = "320018813254"
    10 = "320018813287"
    11 = "320018813302"
    12 = "320018813346"
    13 = "320018813449"
    14 = "320018813471"
    15 = "320018813493"
    16 = "320018813520"
    17 = "320018813541"
    18 = "320018813585"
    19 = "320018813600"
    20 = "320018813850"
    21 = "320018813872"
    22 = "320018813909"
}

# rows where column D (PackageTrackNum) mirrors the new column C value
$rowsWithMirroredD = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in $newTrackNums.Keys) {
    $value = $newTrackNums[$row]

    $cCell = $ws.Cells.Item($row, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $value

    if ($rowsWithMirroredD -contains $row) {
        $dCell = $ws.Cells.Item($row, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $value
    }
}
